$d = $word.ActiveDocument
$sec = $d.Sections(1)
$h = $sec.Headers(2)
$rng = $h.Range
$xml = $rng.WordOpenXML
Write-Output $xml
